$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for every data row (2-520).
# The whole column is bumped by one day: 45204 (2023-10-05) -> 45205 (2023-10-06).
for ($row = 2; $row -le 520; $row++) {
    $ws.Cells.Item($row, 3).Value = 45205
}
